# Sprint 6 burndown update: "4th day and a bit of the 5th too"
# Updates hours-worked data on the "Burndown Chart Sprint6" sheet.
# Downstream formulas (L, N columns; totals in rows 24-26; chart caches)
# recompute automatically from these source values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart Sprint6")

# --- Row 6 ---
$ws.Range("D6").Value = 2
$ws.Range("H6").Value = 0.5
$ws.Range("I6").Value = 0.25

# --- Row 7 ---
$ws.Range("D7").Value = 2
$ws.Range("H7").Value = 0.5

# --- Row 8 ---
$ws.Range("D8").Value = 2
$ws.Range("F8").Value = 0.5
$ws.Range("H8").Value = 0.5

# --- Row 9 ---
$ws.Range("D9").Value = 2
$ws.Range("F9").Value = 0.5
$ws.Range("H9").Value = 0.5

# --- Row 16 ---
$ws.Range("D16").Value = 1

# --- Row 17 ---
$ws.Range("D17").Value = 1
$ws.Range("I17").Value = 0.75

# --- Row 18 ---
$ws.Range("D18").Value = 1

# --- Row 19 ---
$ws.Range("D19").Value = 1

# --- Row 20 ---
$ws.Range("D20").Value = 3
$ws.Range("I20").Value = 1

# --- Row 21 ---
$ws.Range("D21").Value = 3

# --- Row 22 ---
$ws.Range("E22").Value = 0.5
$ws.Range("F22").Value = 0.5

# Put the selection where the author ended up (cell P14) to mirror the
# final UI state after the edits.
$ws.Range("P14").Select()
